$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row 27: Id, ГОСТ, description (added first so new shared strings
# land in the same order as the target workbook)
$ws.Range("A27").Value = 79
$ws.Range("B27").Value = "ГОСТ Р 51613-2000"
$ws.Range("C27").Value = "Трубы напорные из непластифицированного поливинилхлорида. Технические условия"

# New "Версия БД 2.2.2.1" section marker in column C, row 26, centered horizontally
$ws.Range("C26").Value = "Версия БД 2.2.2.1"
$ws.Range("C26").HorizontalAlignment = -4108

# Update the active selection to reflect the new last-used cell
$ws.Range("C26").Select()
